$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 46073
$ws.Range("C3").Value2 = 46073
$ws.Range("C4").Value2 = 46073
$ws.Range("A5").Value2 = "A 16122-2024"
$ws.Range("B5").Value2 = 45406.50112268519
$ws.Range("C5").Value2 = 46073
$ws.Range("G5").Value2 = 15.2
$ws.Range("R5").Value2 = "Oxtungssvamp`r`nMyskmadra"
$ws.Range("S5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/artfynd/A 16122-2024 artfynd.xlsx`", `"A 16122-2024`")"
$ws.Range("T5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/kartor/A 16122-2024 karta.png`", `"A 16122-2024`")"
$ws.Range("V5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/klagomål/A 16122-2024 FSC-klagomål.docx`", `"A 16122-2024`")"
$ws.Range("W5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/klagomålsmail/A 16122-2024 FSC-klagomål mail.docx`", `"A 16122-2024`")"
$ws.Range("X5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/tillsyn/A 16122-2024 tillsynsbegäran.docx`", `"A 16122-2024`")"
$ws.Range("Y5").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/tillsynsmail/A 16122-2024 tillsynsbegäran mail.docx`", `"A 16122-2024`")"
$ws.Range("A6").Value2 = "A 31697-2023"
$ws.Range("B6").Value2 = 45117
$ws.Range("C6").Value2 = 46073
$ws.Range("G6").Value2 = 2.2
$ws.Range("H6").Value2 = 2
$ws.Range("I6").Value2 = 0
$ws.Range("J6").Value2 = 0
$ws.Range("O6").Value2 = 0
$ws.Range("R6").Value2 = "Lövgroda`r`nStörre vattensalamander"
$ws.Range("S6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/artfynd/A 31697-2023 artfynd.xlsx`", `"A 31697-2023`")"
$ws.Range("T6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/kartor/A 31697-2023 karta.png`", `"A 31697-2023`")"
$ws.Range("V6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/klagomål/A 31697-2023 FSC-klagomål.docx`", `"A 31697-2023`")"
$ws.Range("W6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/klagomålsmail/A 31697-2023 FSC-klagomål mail.docx`", `"A 31697-2023`")"
$ws.Range("X6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/tillsyn/A 31697-2023 tillsynsbegäran.docx`", `"A 31697-2023`")"
$ws.Range("Y6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/tillsynsmail/A 31697-2023 tillsynsbegäran mail.docx`", `"A 31697-2023`")"
$ws.Range("C7").Value2 = 46073
$ws.Range("A8").Value2 = "A 15571-2024"
$ws.Range("B8").Value2 = 45401
$ws.Range("C8").Value2 = 46073
$ws.Range("G8").Value2 = 11.1
$ws.Range("H8").Value2 = 0
$ws.Range("I8").Value2 = 1
$ws.Range("J8").Value2 = 1
$ws.Range("O8").Value2 = 1
$ws.Range("R8").Value2 = "Oxtungssvamp`r`nGuldlockmossa"
$ws.Range("S8").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/artfynd/A 15571-2024 artfynd.xlsx`", `"A 15571-2024`")"
$ws.Range("T8").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/kartor/A 15571-2024 karta.png`", `"A 15571-2024`")"
$ws.Range("V8").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/klagomål/A 15571-2024 FSC-klagomål.docx`", `"A 15571-2024`")"
$ws.Range("W8").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/klagomålsmail/A 15571-2024 FSC-klagomål mail.docx`", `"A 15571-2024`")"
$ws.Range("X8").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/tillsyn/A 15571-2024 tillsynsbegäran.docx`", `"A 15571-2024`")"
$ws.Range("Y8").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1263/tillsynsmail/A 15571-2024 tillsynsbegäran mail.docx`", `"A 15571-2024`")"
$ws.Range("C9").Value2 = 46073
$ws.Range("C10").Value2 = 46073
$ws.Range("C11").Value2 = 46073
$ws.Range("C12").Value2 = 46073
$ws.Range("C13").Value2 = 46073
$ws.Range("C14").Value2 = 46073
$ws.Range("C15").Value2 = 46073
$ws.Range("A16").Value2 = "A 15565-2024"
$ws.Range("B16").Value2 = 45401.66103009259
$ws.Range("C16").Value2 = 46073
$ws.Range("G16").Value2 = 7
$ws.Range("A17").Value2 = "A 15582-2024"
$ws.Range("B17").Value2 = 45401.69502314815
$ws.Range("C17").Value2 = 46073
$ws.Range("G17").Value2 = 7.3
$ws.Range("A18").Value2 = "A 25603-2024"
$ws.Range("B18").Value2 = 45463
$ws.Range("C18").Value2 = 46073
$ws.Range("G18").Value2 = 1.8
$ws.Range("A19").Value2 = "A 31701-2023"
$ws.Range("B19").Value2 = 45117
$ws.Range("C19").Value2 = 46073
$ws.Range("G19").Value2 = 0.7
$ws.Range("A20").Value2 = "A 18888-2024"
$ws.Range("B20").Value2 = 45427
$ws.Range("C20").Value2 = 46073
$ws.Range("G20").Value2 = 3.5
$ws.Range("A21").Value2 = "A 31393-2021"
$ws.Range("B21").Value2 = 44368
$ws.Range("C21").Value2 = 46073
$ws.Range("G21").Value2 = 0.6
$ws.Range("A22").Value2 = "A 31693-2023"
$ws.Range("B22").Value2 = 45117
$ws.Range("C22").Value2 = 46073
$ws.Range("G22").Value2 = 4.1
$ws.Range("A23").Value2 = "A 51318-2025"
$ws.Range("B23").Value2 = 45950.38170138889
$ws.Range("C23").Value2 = 46073
$ws.Range("G23").Value2 = 2.8
$ws.Range("A24").Value2 = "A 13344-2023"
$ws.Range("B24").Value2 = 45005
$ws.Range("C24").Value2 = 46073
$ws.Range("G24").Value2 = 1.1
$ws.Range("A25").Value2 = "A 67456-2021"
$ws.Range("B25").Value2 = 44524
$ws.Range("C25").Value2 = 46073
$ws.Range("G25").Value2 = 8.1
$ws.Range("A26").Value2 = "A 2408-2025"
$ws.Range("B26").Value2 = 45674.36113425926
$ws.Range("C26").Value2 = 46073
$ws.Range("G26").Value2 = 1.9
$ws.Range("A27").Value2 = "A 2415-2025"
$ws.Range("B27").Value2 = 45674.379375
$ws.Range("C27").Value2 = 46073
$ws.Range("G27").Value2 = 1.7
$ws.Range("A28").Value2 = "A 16690-2023"
$ws.Range("B28").Value2 = 45030
$ws.Range("C28").Value2 = 46073
$ws.Range("G28").Value2 = 0.9
$ws.Range("A29").Value2 = "A 59432-2022"
$ws.Range("B29").Value2 = 44907
$ws.Range("C29").Value2 = 46073
$ws.Range("G29").Value2 = 5.1
$ws.Range("A30").Value2 = "A 57798-2022"
$ws.Range("C30").Value2 = 46073
$ws.Range("G30").Value2 = 8.9
$ws.Range("A31").Value2 = "A 57803-2022"
$ws.Range("B31").Value2 = 44897
$ws.Range("C31").Value2 = 46073
$ws.Range("A32").Value2 = "A 58350-2025"
$ws.Range("B32").Value2 = 45985.51048611111
$ws.Range("C32").Value2 = 46073
$ws.Range("G32").Value2 = 6.9
$ws.Range("A33").Value2 = "A 29817-2025"
$ws.Range("B33").Value2 = 45825
$ws.Range("C33").Value2 = 46073
$ws.Range("G33").Value2 = 1
$ws.Range("A34").Value2 = "A 15577-2024"
$ws.Range("B34").Value2 = 45401.68829861111
$ws.Range("C34").Value2 = 46073
$ws.Range("G34").Value2 = 1.9
$ws.Range("A35").Value2 = "A 12953-2025"
$ws.Range("B35").Value2 = 45734
$ws.Range("C35").Value2 = 46073
$ws.Range("G35").Value2 = 1.4
$ws.Range("A36").Value2 = "A 19190-2023"
$ws.Range("B36").Value2 = 45048
$ws.Range("C36").Value2 = 46073
$ws.Range("G36").Value2 = 0.5
$ws.Range("A37").Value2 = "A 16199-2024"
$ws.Range("B37").Value2 = 45406
$ws.Range("C37").Value2 = 46073
$ws.Range("G37").Value2 = 13.5
$ws.Range("A38").Value2 = "A 54580-2023"
$ws.Range("B38").Value2 = 45233
$ws.Range("C38").Value2 = 46073
$ws.Range("G38").Value2 = 0.5
$ws.Range("A39").Value2 = "A 64060-2025"
$ws.Range("B39").Value2 = 46021.6172337963
$ws.Range("C39").Value2 = 46073
$ws.Range("G39").Value2 = 6.4
$ws.Range("A40").Value2 = "A 64055-2025"
$ws.Range("B40").Value2 = 46021
$ws.Range("C40").Value2 = 46073
$ws.Range("G40").Value2 = 0.9
$ws.Range("A41").Value2 = "A 64049-2025"
$ws.Range("B41").Value2 = 46021
$ws.Range("C41").Value2 = 46073
$ws.Range("G41").Value2 = 0.8
$ws.Range("A42").Value2 = "A 64051-2025"
$ws.Range("B42").Value2 = 46021
$ws.Range("C42").Value2 = 46073
$ws.Range("G42").Value2 = 1
$ws.Range("A43").Value2 = "A 64058-2025"
$ws.Range("B43").Value2 = 46021
$ws.Range("C43").Value2 = 46073
$ws.Range("G43").Value2 = 1.4
$ws.Range("A44").Value2 = "A 9570-2026"
$ws.Range("B44").Value2 = 46070
$ws.Range("C44").Value2 = 46073
$ws.Range("G44").Value2 = 6.6
$ws.Range("C45").Value2 = 46073
$ws.Range("A46").Value2 = "A 54284-2023"
$ws.Range("B46").Value2 = 45232
$ws.Range("C46").Value2 = 46073
$ws.Range("G46").Value2 = 2.2
$ws.Range("A47").Value2 = "A 56202-2023"
$ws.Range("B47").Value2 = 45240
$ws.Range("C47").Value2 = 46073
$ws.Range("G47").Value2 = 2.8
$ws.Range("A48").Value2 = "A 38011-2023"
$ws.Range("B48").Value2 = 45160
$ws.Range("C48").Value2 = 46073
$ws.Range("G48").Value2 = 2.8
$ws.Range("A49").Value2 = "A 54280-2023"
$ws.Range("B49").Value2 = 45232.69518518518
$ws.Range("C49").Value2 = 46073
$ws.Range("G49").Value2 = 1.4
$ws.Range("A50").Value2 = "A 15277-2024"
$ws.Range("B50").Value2 = 45400
$ws.Range("C50").Value2 = 46073
$ws.Range("G50").Value2 = 1.3
$ws.Range("A51").Value2 = "A 27113-2022"
$ws.Range("B51").Value2 = 44741
$ws.Range("C51").Value2 = 46073
$ws.Range("G51").Value2 = 8.9
$ws.Range("A52").Value2 = "A 29030-2023"
$ws.Range("B52").Value2 = 45104
$ws.Range("C52").Value2 = 46073
$ws.Range("A53").Value2 = "A 54282-2023"
$ws.Range("B53").Value2 = 45232.69699074074
$ws.Range("C53").Value2 = 46073
$ws.Range("G53").Value2 = 1.8
$ws.Range("A54").Value2 = "A 9277-2025"
$ws.Range("B54").Value2 = 45714.63053240741
$ws.Range("C54").Value2 = 46073
$ws.Range("G54").Value2 = 1.7
